$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column AB
$ws.Range("AB1").Value = "12-jul"

# New data values for column AB (rows 2-18)
$values = @(
    0,
    11.798624784142465,
    20.164695214552214,
    24.755601365228237,
    0,
    11.534424368550971,
    10.125520551384172,
    24.08038975539549,
    21.676717950743591,
    13.108048188524972,
    0,
    12.929436441383448,
    0,
    0,
    20.531615949268769,
    0,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 28).Value = $values[$i]
}

# Update selection to match the new active range
$ws.Range("AB2:AB18").Select()
